# Update the Bento multi-filter query cells on the "startup" sheet.
# Write order matters for shared-string table compaction/append behavior,
# so B4 (File query), then B3 (Sample query), then B2 (Case query) are
# written last-to-first relative to their row position, matching the
# target shared-string layout produced by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "MATCH (f:file)`nMATCH (f)-[:file_of_sample]->(samp)`nMATCH (samp)-[:sample_of_study_subject]->(ss)`nMATCH (ss)-[:study_subject_of_study]->(s)`nMATCH (ss)<-[:sf_of_study_subject]-(sf)`nMATCH (s)-[:study_of_program]->(p)`nMATCH (d)-[:diagnosis_of_study_subject]->(ss)`nMATCH (tp)-[:tp_of_diagnosis]->(d)`nWHERE s.study_acronym IN [`"C`"] and ss.disease_subtype IN [`"Papillary Carcinoma`"] and d.tumor_size_group In [`"<=1`"] and d.pr_status In [`"Positive`"] and tp.endocrine_therapy_type IN [`"Tam`"]`nRETURN  f.file_name AS ``File Name``,`n    head(labels(samp)) AS ``Association``,`n    f.file_description AS ``Description``,`n    f.file_format AS ``File Format``,`n    f.file_size AS ``Size``,`n    p.program_acronym AS ``Program Code``,`n    s.study_acronym AS ``Arm``,`n    ss.study_subject_id AS ``Case ID``,`n    samp.sample_id AS ``Sample ID``,`n    ss.disease_subtype as ``Diagnosis```n order By f.file_name ASC LIMIT 100"
$ws.Range("B3").Value = "MATCH (ss:study_subject)`nWITH COLLECT(ss.study_subject_id) AS all_subjects`nMATCH (samp:sample)`nMATCH (samp)-[:sample_of_study_subject]->(ss)`nMATCH (ss)<-[:sf_of_study_subject]-(sf)`nMATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)`nMATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)`nMATCH (ss)<-[:diagnosis_of_study_subject]-(d)`nMATCH (d)<-[:tp_of_diagnosis]-(tp)`nWHERE s.study_acronym IN [`"C`"] and ss.disease_subtype IN [`"Papillary Carcinoma`"] and d.tumor_size_group In [`"<=1`"] and d.pr_status In [`"Positive`"] and tp.endocrine_therapy_type IN [`"Tam`"]`nWITH`n    distinct lp,`n    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,`n    collect(distinct f.file_id) AS files,`n    samp, ss, s, p, all_subjects`nRETURN`n samp.sample_id AS ``Sample ID``,`n            ss.study_subject_id AS ``Case ID``,`n            p.program_acronym AS ``Program Code``,`n            s.study_acronym AS ``Arm``,`n            ss.disease_subtype AS ``Diagnosis``,`n            samp.tissue_type AS ``Tissue Type``,`n            samp.composition AS ``Tissue Composition``,`n            samp.sample_anatomic_site AS ``Sample Anatomic Site``,`n            samp.method_of_sample_procurement AS ``Sample Procurement Method``,`n            lp.test_name as Platform`n order By samp.sample_id ASC LIMIT 100"
$ws.Range("B2").Value = "MATCH (ss:study_subject)`nMATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)`nWITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files`nMATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)`nMATCH (ss)<-[:sf_of_study_subject]-(sf)`nMATCH (ss)<-[:diagnosis_of_study_subject]-(d)`nMATCH (d)<-[:tp_of_diagnosis]-(tp)`nMATCH (ss)<-[:demographic_of_study_subject]-(demo)`nWHERE s.study_acronym IN [`"C`"] and ss.disease_subtype IN [`"Papillary Carcinoma`"] and d.tumor_size_group In [`"<=1`"] and d.pr_status In [`"Positive`"] and tp.endocrine_therapy_type IN [`"Tam`"]`nreturn ss.study_subject_id as ``Case ID``,`n       p.program_acronym as ``Program Code``,`n        p.program_id as Program_ID,`n       s.study_acronym as ``Arm``,`n       ss.disease_subtype as ``Diagnosis``,`n       sf.grouped_recurrence_score AS ``Recurrence Score``,`n       d.tumor_size_group AS ``tumor_size``,`n       d.er_status AS ``ER Status``,`n       d.pr_status AS ``PR Status``,`n       demo.age_at_index AS ``Age (years)``,`ndemo.survival_time AS ``Survival (days)```n order By ss.study_subject_id ASC LIMIT 100"

# Row heights grew along with the appended "order By ... LIMIT 100" text.
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 409.6
$ws.Rows.Item(4).RowHeight = 316.8

# Selection moved from D2 to B2.
$ws.Range("B2").Select()
